$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (24) to the Question_List sheet for "Transform to Sum Tree"
$ws.Range("A24").Value = "GFG"
$ws.Range("B24").Value = "GFG"
$ws.Range("C24").Value = "Transform to Sum Tree"

# Match styling used by similar rows (e.g. row 18/20): left-aligned, wrap text
$ws.Range("A24").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B24:C24").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B24:C24").VerticalAlignment = -4160     # xlTop
$ws.Range("B24:C24").WrapText = $true

# Update the selected cell to reflect where the user ended up
$ws.Range("D27").Select()
